$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(1, 0.6834618577321858),
    @(2, 0.614852389721919),
    @(3, 0.5772299978691616),
    @(5, 0.5644555014650496),
    @(6, 0.5763062362569621),
    @(7, 0.5759841008387409),
    @(8, 0.5759317693843864),
    @(10, 0.5810875563963549)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
